# Auto-generated Excel COM script to apply market price data updates
# across multiple worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

# ALC!row86
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 3063.5
$ws.Range("I86").Value = 1405.25
$ws.Range("J86").Value = 4058.45
$ws.Range("K86").Value = 1405.25
$ws.Range("L86").Value = 4058.45
$ws.Range("M86").Value = -282.25
$ws.Range("N86").Value = -6304.45

# ALC!row89
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 3063.5
$ws.Range("I89").Value = 1405.25
$ws.Range("J89").Value = 4058.45
$ws.Range("K89").Value = 7026.25
$ws.Range("L89").Value = 20292.25
$ws.Range("M89").Value = -1410.25
$ws.Range("N89").Value = -31524.25

# ALC!row113
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 2339.1924
$ws.Range("I113").Value = 1623.9231
$ws.Range("J113").Value = 3054.4614
$ws.Range("K113").Value = 1623.9231
$ws.Range("L113").Value = 3054.4614
$ws.Range("M113").Value = 1630.0769
$ws.Range("N113").Value = -9562.4614

# ALC!row114
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H114").Value = 48000
$ws.Range("J114").Value = 48000
$ws.Range("L114").Value = 48000
$ws.Range("N114").Value = -56678

# ALC!row125
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 8477.706
$ws.Range("I125").Value = 15964
$ws.Range("J125").Value = 1823.2222
$ws.Range("K125").Value = 143676
$ws.Range("L125").Value = 16408.9998
$ws.Range("M125").Value = -141216
$ws.Range("N125").Value = -21328.9998

# ARM!row32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4380.9897
$ws.Range("I32").Value = 3539.875
$ws.Range("K32").Value = 3539.875
$ws.Range("M32").Value = -3252.875

# ARM!row34
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H34").Value = 8000
$ws.Range("I34").Value = 8000
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 8000
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -7729
$ws.Range("N34").ClearContents()

# ARM!row74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2777.3674
$ws.Range("I74").Value = 638.7931
$ws.Range("J74").Value = 5878.3
$ws.Range("K74").Value = 638.7931
$ws.Range("L74").Value = 5878.3
$ws.Range("M74").Value = 235.2069
$ws.Range("N74").Value = -7626.3

# ARM!row77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 2777.3674
$ws.Range("I77").Value = 638.7931
$ws.Range("J77").Value = 5878.3
$ws.Range("K77").Value = 3193.9655
$ws.Range("L77").Value = 29391.5
$ws.Range("M77").Value = 1174.0345
$ws.Range("N77").Value = -38127.5

# ARM!row113
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H113").Value = 48000
$ws.Range("J113").Value = 48000
$ws.Range("L113").Value = 48000
$ws.Range("N113").Value = -56678

# ARM!row114
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H114").Value = 28350
$ws.Range("J114").Value = 28350
$ws.Range("L114").Value = 28350
$ws.Range("N114").Value = -37028

# ARM!row132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 5091.314
$ws.Range("I132").Value = 3579.923
$ws.Range("J132").Value = 10003.333
$ws.Range("K132").Value = 10739.769
$ws.Range("L132").Value = 30009.999
$ws.Range("M132").Value = -8209.769
$ws.Range("N132").Value = -35069.999

# ARM!row139
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H139").Value = 80000
$ws.Range("J139").Value = 80000
$ws.Range("L139").Value = 80000
$ws.Range("N139").Value = -90280

# ARM!row140
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H140").Value = 50000
$ws.Range("J140").Value = 50000
$ws.Range("L140").Value = 50000
$ws.Range("N140").Value = -60360

# ARM!row141
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H141").Value = 50000
$ws.Range("J141").Value = 50000
$ws.Range("L141").Value = 50000
$ws.Range("N141").Value = -60360

# BSM!row114
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H114").Value = 47500
$ws.Range("J114").Value = 47500
$ws.Range("L114").Value = 47500
$ws.Range("N114").Value = -56178

# CRP!row58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 801.8644
$ws.Range("I58").Value = 654.63416
$ws.Range("K58").Value = 654.63416
$ws.Range("M58").Value = -451.63416

# CRP!row62
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 3971.795
$ws.Range("I62").Value = 4761.087
$ws.Range("J62").Value = 2837.1875
$ws.Range("K62").Value = 4761.087
$ws.Range("L62").Value = 2837.1875
$ws.Range("M62").Value = -4137.087
$ws.Range("N62").Value = -4085.1875

# CRP!row65
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 3971.795
$ws.Range("I65").Value = 4761.087
$ws.Range("J65").Value = 2837.1875
$ws.Range("K65").Value = 23805.435
$ws.Range("L65").Value = 14185.9375
$ws.Range("M65").Value = -20685.435
$ws.Range("N65").Value = -20425.9375

# CRP!row99
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 3735.4285
$ws.Range("I99").Value = 3524.6667
$ws.Range("K99").Value = 3524.6667
$ws.Range("M99").Value = -2026.6667

# CRP!row126
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 3735.4285
$ws.Range("I126").Value = 3524.6667
$ws.Range("K126").Value = 10574.0001
$ws.Range("M126").Value = -8104.000100000001

# CRP!row134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 3308.1836
$ws.Range("I134").Value = 3448.6
$ws.Range("J134").Value = 1728.5
$ws.Range("K134").Value = 10345.8
$ws.Range("L134").Value = 5185.5
$ws.Range("M134").Value = -7810.799999999999
$ws.Range("N134").Value = -10255.5

# CRP!row136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 801.8644
$ws.Range("I136").Value = 654.63416
$ws.Range("K136").Value = 1963.90248
$ws.Range("M136").Value = 586.0975200000003

# CUL!row113
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 656821.3
$ws.Range("I113").Value = 447.2857
$ws.Range("J113").Value = 1575744.9
$ws.Range("K113").Value = 1341.8571
$ws.Range("L113").Value = 4727234.699999999
$ws.Range("M113").Value = 828.1428999999998
$ws.Range("N113").Value = -4731574.699999999

# CUL!row114
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value = 823.6667
$ws.Range("J114").Value = 1932
$ws.Range("L114").Value = 5796
$ws.Range("N114").Value = -12304

# CUL!row122
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 794.2833000000001
$ws.Range("I122").Value = 364.875
$ws.Range("J122").Value = 950.4318
$ws.Range("K122").Value = 3283.875
$ws.Range("L122").Value = 8553.886199999999
$ws.Range("M122").Value = -833.875
$ws.Range("N122").Value = -13453.8862

# GSM!row113
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1655
$ws.Range("I113").Value = 1666.6666
$ws.Range("J113").Value = 1650
$ws.Range("K113").Value = 1666.6666
$ws.Range("L113").Value = 1650
$ws.Range("M113").Value = 503.3334
$ws.Range("N113").Value = -5990

# GSM!row114
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H114").Value = 41500
$ws.Range("J114").Value = 41500
$ws.Range("L114").Value = 41500
$ws.Range("N114").Value = -50178

# LTW!row7
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1742
$ws.Range("I7").Value = 1331.8462
$ws.Range("J7").Value = 3075
$ws.Range("K7").Value = 1331.8462
$ws.Range("L7").Value = 3075
$ws.Range("M7").Value = -1219.8462
$ws.Range("N7").Value = -3299

# LTW!row61
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1363.0571
$ws.Range("I61").Value = 1051.3334
$ws.Range("J61").Value = 1596.85
$ws.Range("K61").Value = 1051.3334
$ws.Range("L61").Value = 1596.85
$ws.Range("M61").Value = -849.3334
$ws.Range("N61").Value = -2000.85

# LTW!row68
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1846.8518
$ws.Range("I68").Value = 1754.2354
$ws.Range("J68").Value = 2004.3
$ws.Range("K68").Value = 1754.2354
$ws.Range("L68").Value = 2004.3
$ws.Range("M68").Value = -1005.2354
$ws.Range("N68").Value = -3502.3

# LTW!row71
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 1846.8518
$ws.Range("I71").Value = 1754.2354
$ws.Range("J71").Value = 2004.3
$ws.Range("K71").Value = 8771.177
$ws.Range("L71").Value = 10021.5
$ws.Range("M71").Value = -5027.177
$ws.Range("N71").Value = -17509.5

# LTW!row113
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 1363.0571
$ws.Range("I113").Value = 1051.3334
$ws.Range("J113").Value = 1596.85
$ws.Range("K113").Value = 1051.3334
$ws.Range("L113").Value = 1596.85
$ws.Range("M113").Value = 1118.6666
$ws.Range("N113").Value = -5936.85

# LTW!row114
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()

# LTW!row115
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H115").Value = 40000
$ws.Range("J115").Value = 40000
$ws.Range("L115").Value = 40000
$ws.Range("N115").Value = -42350

# LTW!row126
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 1742
$ws.Range("I126").Value = 1331.8462
$ws.Range("J126").Value = 3075
$ws.Range("K126").Value = 3995.5386
$ws.Range("L126").Value = 9225
$ws.Range("M126").Value = -1525.5386
$ws.Range("N126").Value = -14165

# LTW!row132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3097.0127
$ws.Range("I132").Value = 3152.3015
$ws.Range("J132").Value = 2879.3125
$ws.Range("K132").Value = 9456.904500000001
$ws.Range("L132").Value = 8637.9375
$ws.Range("M132").Value = -6926.904500000001
$ws.Range("N132").Value = -13697.9375

# LTW!row136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 2872.9092
$ws.Range("I136").Value = 907.5217
$ws.Range("J136").Value = 12918.223
$ws.Range("K136").Value = 2722.5651
$ws.Range("L136").Value = 38754.669
$ws.Range("M136").Value = -172.5650999999998
$ws.Range("N136").Value = -43854.669

# WVR!row113
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 48062.094
$ws.Range("I113").Value = 77121.84
$ws.Range("J113").Value = 840
$ws.Range("K113").Value = 231365.52
$ws.Range("L113").Value = 2520
$ws.Range("M113").Value = -229195.52
$ws.Range("N113").Value = -6860

# WVR!row114
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H114").Value = 30000
$ws.Range("J114").Value = 30000
$ws.Range("L114").Value = 30000
$ws.Range("N114").Value = -38678

# WVR!row126
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 717
$ws.Range("I126").Value = 487.5
$ws.Range("J126").Value = 1307.1428
$ws.Range("K126").Value = 1462.5
$ws.Range("L126").Value = 3921.4284
$ws.Range("M126").Value = 1007.5
$ws.Range("N126").Value = -8861.428400000001

# WVR!row136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 2862.4805
$ws.Range("I136").Value = 3801.5715
$ws.Range("J136").Value = 1735.5714
$ws.Range("K136").Value = 11404.7145
$ws.Range("L136").Value = 5206.7142
$ws.Range("M136").Value = -8854.7145
$ws.Range("N136").Value = -10306.7142

Write-Output "Applied all market price updates."